# This script reproduces the diff: two new data rows (for a new reporting date)
# are inserted right before the existing row 1058, pushing all subsequent rows
# down by two (old row N -> new row N+2). The worksheet dimension grows from
# A1:R1190 to A1:R1192.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before row 1058 (do it twice so both inherit formatting
# from the row that was originally at 1058, in particular the date style on D).
$ws.Rows.Item(1058).Insert()
$ws.Rows.Item(1058).Insert()

# Constant values shared by every data row in this sheet.
$mercadoId = 8
$mercado   = "Terminal La Palmera de La Serena"
$region    = "Coquimbo"
$codreg    = 4
$catId     = 100112023
$categoria = "Brócoli"
$variedad  = "Sin especificar"
$unidad    = "`$/unidad"
$origen    = "Provincia del Elquí"
$kgUnid    = 1
$clasif    = "Hortaliza"

# New row 1058: "Primera" quality entry for the new date (serial 45142).
$r = 1058
$ws.Cells.Item($r,1).Value  = $mercadoId
$ws.Cells.Item($r,2).Value  = $mercado
$ws.Cells.Item($r,3).Value  = $region
$ws.Cells.Item($r,4).Value  = 45142
$ws.Cells.Item($r,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($r,5).Value  = $codreg
$ws.Cells.Item($r,6).Value  = $catId
$ws.Cells.Item($r,7).Value  = $categoria
$ws.Cells.Item($r,8).Value  = $variedad
$ws.Cells.Item($r,9).Value  = "Primera"
$ws.Cells.Item($r,10).Value = 2040
$ws.Cells.Item($r,11).Value = 700
$ws.Cells.Item($r,12).Value = 800
$ws.Cells.Item($r,13).Value = 750
$ws.Cells.Item($r,14).Value = $unidad
$ws.Cells.Item($r,15).Value = $origen
$ws.Cells.Item($r,16).Value = 750
$ws.Cells.Item($r,17).Value = $kgUnid
$ws.Cells.Item($r,18).Value = $clasif

# New row 1059: "Segunda" quality entry for the same new date.
$r = 1059
$ws.Cells.Item($r,1).Value  = $mercadoId
$ws.Cells.Item($r,2).Value  = $mercado
$ws.Cells.Item($r,3).Value  = $region
$ws.Cells.Item($r,4).Value  = 45142
$ws.Cells.Item($r,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($r,5).Value  = $codreg
$ws.Cells.Item($r,6).Value  = $catId
$ws.Cells.Item($r,7).Value  = $categoria
$ws.Cells.Item($r,8).Value  = $variedad
$ws.Cells.Item($r,9).Value  = "Segunda"
$ws.Cells.Item($r,10).Value = 1300
$ws.Cells.Item($r,11).Value = 500
$ws.Cells.Item($r,12).Value = 600
$ws.Cells.Item($r,13).Value = 550
$ws.Cells.Item($r,14).Value = $unidad
$ws.Cells.Item($r,15).Value = $origen
$ws.Cells.Item($r,16).Value = 550
$ws.Cells.Item($r,17).Value = $kgUnid
$ws.Cells.Item($r,18).Value = $clasif
